$wb = $excel.ActiveWorkbook

# Delete rows 2-4 (the "Accrual" transaction rows) from the Transactions sheet;
# rows 5 & 6 shift up to become the new rows 2 & 3.
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Rows("2:4").Delete()

# Leave behind a selection on the Repayment schedule sheet.
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$null = $wsRepay.Range("C4").Select()

# Select the region that used to hold the deleted rows, and make Transactions
# the active sheet/tab (last Select() wins for the active sheet).
$null = $wsTrans.Range("A2:XFD4").Select()
